# Append the 2025-11-04 allocation row (row 64) to the sheet, matching the
# existing rows: a date stored as plain text in column A and the two
# allocation fractions as numbers in columns B and C.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading apostrophe forces the date-looking string to be entered as text
# instead of being auto-converted to a date serial number (mirrors how the
# other rows store their Date column as plain text).
$ws.Range("A64").Value = "'11/04/2025"

# Re-apply the plain "no explicit format" style used by the other data rows
# so the new cell doesn't pick up the quote-prefix formatting that typing a
# leading apostrophe would otherwise leave behind.
$ws.Range("A64").Style = $ws.Range("A63").Style

$ws.Range("B64").Value = 0.2154778411495352
$ws.Range("C64").Value = 0.7845221588504648
